# Extent reports updated to show categories
# Remove the "Create Enrollment" test case row from Sheet2 (it was row 8),
# which shifts the subsequent rows (My Profile, CSR Manage Users, CSR View
# Payments) up by one, and updates the sheet's selection/dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Delete the entire row 8 ("UPA_Regression" | "Create Enrollment" |
# "test.java.TestCreateEnrollment"), shifting rows below it up.
$ws.Rows.Item(8).Delete()

# Update the active selection to match the new last data row (C10).
$ws.Range("C10").Select()
